$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the formula in A2 with a literal value (0)
$ws.Range("A2").Value = 0

# Update the selected cell to B2
$ws.Range("B2").Select()
